$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match the style of the date column used by the other rows (numFmtId 49 / text format)
# before writing the value, so the date-like string is kept as literal text
# instead of being auto-converted to a date serial number.
$ws.Range("A37").Style = $ws.Range("A36").Style
$ws.Range("A37").NumberFormat = "@"

# Add the new row of data for 2020-06-16 (report updated 18 Junio, no municipality breakdown)
$ws.Range("A37").Value = "2020-06-16"
$ws.Range("B37").Value = "Andalucía"
$ws.Range("C37").Value = 49
$ws.Range("D37").Value = 14

# Update selection to match the recorded cursor position after the edit
$ws.Range("C38").Select()
